$wb = $excel.ActiveWorkbook

# Update profit-calculation columns (H-N) on the Leve profit sheets.
# Values come from the latest Jenova market-price refresh.

# --- ALC sheet, row 40 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4975.05
$ws.Range("I40").Value = 2346.2307
$ws.Range("J40").Value = 9857.143
$ws.Range("K40").Value = 2346.2307
$ws.Range("L40").Value = 9857.143
$ws.Range("M40").Value = -2171.2307
$ws.Range("N40").Value = -10207.143

# --- ALC sheet, row 112 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3336.88
$ws.Range("I112").Value = 1777
$ws.Range("J112").Value = 3401.875
$ws.Range("K112").Value = 5331
$ws.Range("L112").Value = 10205.625
$ws.Range("M112").Value = -4223
$ws.Range("N112").Value = -12421.625

# --- ALC sheet, row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3642.7908
$ws.Range("I137").Value = 2713.182
$ws.Range("J137").Value = 4616.6665
$ws.Range("K137").Value = 8139.545999999999
$ws.Range("L137").Value = 13849.9995
$ws.Range("M137").Value = -5589.545999999999
$ws.Range("N137").Value = -18949.9995

# --- ALC sheet, row 141 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4798.5713
$ws.Range("I141").Value = 5098.3335
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 15295.0005
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -10115.0005
$ws.Range("N141").Value = -19360

# --- ARM sheet, row 27 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 9750
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 9750
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 9750
$ws.Range("N27").Value = -10118

# --- ARM sheet, row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3393.868
$ws.Range("I32").Value = 2742.647
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 2742.647
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -2455.647
$ws.Range("N32").Value = -20574

# --- ARM sheet, row 61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3728.5833
$ws.Range("I61").Value = 1969.7646
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 1969.7646
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -1757.7646
$ws.Range("N61").Value = -8424

# --- ARM sheet, row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 297935.1
$ws.Range("I74").Value = 402901.03
$ws.Range("J74").Value = 6363
$ws.Range("K74").Value = 402901.03
$ws.Range("L74").Value = 6363
$ws.Range("M74").Value = -402027.03
$ws.Range("N74").Value = -8111

# --- ARM sheet, row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 297935.1
$ws.Range("I77").Value = 402901.03
$ws.Range("J77").Value = 6363
$ws.Range("K77").Value = 2014505.15
$ws.Range("L77").Value = 31815
$ws.Range("M77").Value = -2010137.15
$ws.Range("N77").Value = -40551

# --- ARM sheet, row 135 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()

# --- ARM sheet, row 136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3728.5833
$ws.Range("I136").Value = 1969.7646
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 5909.293799999999
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -3359.293799999999
$ws.Range("N136").Value = -29100

# --- BSM sheet, row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17431.455
$ws.Range("I134").Value = 1520.9445
$ws.Range("J134").Value = 89028.75
$ws.Range("K134").Value = 4562.833500000001
$ws.Range("L134").Value = 267086.25
$ws.Range("M134").Value = -2027.833500000001
$ws.Range("N134").Value = -272156.25

# --- CRP sheet, row 19 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4999.2856
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4999.2856
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 4999.2856
$ws.Range("N19").Value = -5339.2856

# --- CRP sheet, row 24 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 4999.2856
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 4999.2856
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 4999.2856
$ws.Range("N24").Value = -5339.2856

# --- CRP sheet, row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4212.6665
$ws.Range("I31").Value = 2115
$ws.Range("J31").Value = 5111.6665
$ws.Range("K31").Value = 2115
$ws.Range("L31").Value = 5111.6665
$ws.Range("M31").Value = -1820
$ws.Range("N31").Value = -5701.6665

# --- CRP sheet, row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4212.6665
$ws.Range("I34").Value = 2115
$ws.Range("J34").Value = 5111.6665
$ws.Range("K34").Value = 2115
$ws.Range("L34").Value = 5111.6665
$ws.Range("M34").Value = -1913
$ws.Range("N34").Value = -5515.6665

# --- CRP sheet, row 58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 206824.4
$ws.Range("I58").Value = 334965.47
$ws.Range("J58").Value = 4496.421
$ws.Range("K58").Value = 334965.47
$ws.Range("L58").Value = 4496.421
$ws.Range("M58").Value = -334762.47
$ws.Range("N58").Value = -4902.421

# --- CRP sheet, row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4976.231
$ws.Range("I132").Value = 5166.4116
$ws.Range("J132").Value = 3683
$ws.Range("K132").Value = 15499.2348
$ws.Range("L132").Value = 11049
$ws.Range("M132").Value = -12969.2348
$ws.Range("N132").Value = -16109

# --- CRP sheet, row 136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 206824.4
$ws.Range("I136").Value = 334965.47
$ws.Range("J136").Value = 4496.421
$ws.Range("K136").Value = 1004896.41
$ws.Range("L136").Value = 13489.263
$ws.Range("M136").Value = -1002346.41
$ws.Range("N136").Value = -18589.263

# --- CUL sheet, row 11 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 125368.875
$ws.Range("I11").Value = 133702.94
$ws.Range("J11").Value = 358
$ws.Range("K11").Value = 401108.82
$ws.Range("L11").Value = 1074
$ws.Range("M11").Value = -400968.82
$ws.Range("N11").Value = -1354

# --- CUL sheet, row 68 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3066.8333
$ws.Range("I68").Value = 2266.9167
$ws.Range("J68").Value = 4666.6665
$ws.Range("K68").Value = 6800.750100000001
$ws.Range("L68").Value = 13999.9995
$ws.Range("M68").Value = -5989.750100000001
$ws.Range("N68").Value = -15621.9995

# --- CUL sheet, row 71 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3066.8333
$ws.Range("I71").Value = 2266.9167
$ws.Range("J71").Value = 4666.6665
$ws.Range("K71").Value = 20402.2503
$ws.Range("L71").Value = 41999.9985
$ws.Range("M71").Value = -16346.2503
$ws.Range("N71").Value = -50111.9985

# --- LTW sheet, row 19 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 166671330
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 166671330
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 166671330
$ws.Range("N19").Value = -166671670

# --- LTW sheet, row 32 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 60
$ws.Range("I32").Value = 60
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 60
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 257

# --- LTW sheet, row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4858.7915
$ws.Range("I132").Value = 3753.8
$ws.Range("J132").Value = 6700.4443
$ws.Range("K132").Value = 11261.4
$ws.Range("L132").Value = 20101.3329
$ws.Range("M132").Value = -8731.400000000001
$ws.Range("N132").Value = -25161.3329

# --- WVR sheet, row 54 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 60000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 60000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 60000
$ws.Range("N54").Value = -61040

# --- WVR sheet, row 107 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 28441.432
$ws.Range("I107").Value = 45191.523
$ws.Range("J107").Value = 923.4286
$ws.Range("K107").Value = 135574.569
$ws.Range("L107").Value = 2770.2858
$ws.Range("M107").Value = -133654.569
$ws.Range("N107").Value = -6610.2858

Write-Output "Updated Leve profit values across ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets."